$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- A1: "phone" keeps its value, gains a thin border on all sides ---
$ws.Range("A1").Borders.LineStyle = 1

# --- A2: "082164886204" keeps its value, loses its special font, gains
#         Text number format ("@") and a thin border ---
$ws.Range("A2").Style = "Normal"
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Borders.LineStyle = 1

# --- New text values populated in the same order the original author
#     typed them in (this controls the shared-string table order) ---

# A8 : "abcdefghi"
$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = "abcdefghi"
$ws.Range("A8").Borders.LineStyle = 1

# A9 : "null"
$ws.Range("A9").NumberFormat = "@"
$ws.Range("A9").Value = "null"
$ws.Range("A9").Borders.LineStyle = 1

# A3 : "6282164886204"
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "6282164886204"
$ws.Range("A3").Borders.LineStyle = 1

# A4 : "+6282164886204"
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "+6282164886204"
$ws.Range("A4").Borders.LineStyle = 1

# A5 : "08216488"
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "08216488"
$ws.Range("A5").Borders.LineStyle = 1

# A6 : "08216488620444"
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "08216488620444"
$ws.Range("A6").Borders.LineStyle = 1

# A7 : "082111112222" - also left-aligned
$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = "082111112222"
$ws.Range("A7").Borders.LineStyle = 1
$ws.Range("A7").HorizontalAlignment = -4131

# A10 : "082164886204" (reuses the shared string already used by A2)
$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = "082164886204"
$ws.Range("A10").Borders.LineStyle = 1

# --- Selection / active cell moves to C8 ---
$ws.Range("C8").Select()

Write-Host "db assertion rows updated"
